$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K1").Value = "Usage_Date"
$ws.Range("L1").Value = "Usage_Count"
$ws.Range("K1:L1").Select()
